$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.315.47"
$ws.Range("E2").Value = "  -1.04%  "

$ws.Range("D3").Value = "1.839.28"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6247"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.03%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07371"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.54%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2882"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.24%  "

$ws.Range("E10").Value = "  -1.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07727"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.32%  "

$ws.Range("D12").Value = "1.839.51"
$ws.Range("E12").Value = "  -0.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.948"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001053"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6620"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.242"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.67%  "

$ws.Range("D18").Value = "29.282.41"
$ws.Range("E18").Value = "  -1.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "233.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.289"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.418"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1335"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07114"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.486"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.480"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.46%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.022"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.30%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.027"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.149"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.808"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6945"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.583"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01825"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.783"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.14%  "

$ws.Range("D39").Value = "1.232.78"
$ws.Range("E39").Value = "  -2.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.783"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9493"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.13%  "

$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.32%  "

$ws.Range("D44").Value = "1.985.83"
$ws.Range("E44").Value = "  -2.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000117"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.932"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.75%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.675"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.918"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1129"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3867"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.28%  "
